$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: A13 was mistakenly logged as text "71277628" -- fix it in place to
# be the numeric phone id, matching the other rows.
$ws.Cells.Item(13, 1).Value = 71277628

# Redeem points 71277628 76.0 -- append the new redemption log entry as row 14.
# Phone / timestamp are kept as text (leading apostrophe forces Excel to store
# them as literal strings rather than re-interpreting as numbers/dates), and
# points is numeric, matching the shape of the existing log rows.
$ws.Cells.Item(14, 1).Value = "'71277628"
$ws.Cells.Item(14, 2).Value = 76
$ws.Cells.Item(14, 3).Value = "'2025-08-18T16:54:54"
